# Penalty Reward System (unfinished) - shift forecast weeks forward by one
# and update the recalculated MyForecast / Summary figures.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$weekDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

$myForecast = @(
    20,
    21,
    21,
    20,
    19,
    18,
    18,
    17,
    16,
    15,
    15,
    14,
    13,
    12,
    12,
    11
)

for ($i = 0; $i -lt $weekDates.Length; $i++) {
    $row = $i + 2
    $dateCell = $ws1.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $weekDates[$i]
    $ws1.Cells.Item($row, 4).Value = $myForecast[$i]
}

# --- Sheet 2: "Summary" ---
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "2023-01-01 to 2025-01-05"

$ws2.Range("B4").NumberFormat = "@"
$ws2.Range("B4").Value = "78"

$ws2.Range("B8").NumberFormat = "@"
$ws2.Range("B8").Value = "3318 units"

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "262"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "154"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "82"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "21"

$ws2.Range("B13").NumberFormat = "@"
$ws2.Range("B13").Value = "2025-01-19"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "11"

$ws2.Range("B15").NumberFormat = "@"
$ws2.Range("B15").Value = "2025-04-27"
